# Update countries & provincias Spain
# Applies the data refresh captured in the commit: updated case counts for
# several countries, an updated "last refreshed" timestamp, and the
# Islas Malvinas / Groenlandia row-order swap in the shared country list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "datos actualizados" timestamp (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 06:50"

# --- Swap the Islas Malvinas / Groenlandia rows (order changed upstream) ---
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Brasil (row 5): Casos activos / Recuperados ---
$ws.Range("D5").Value = 1366775
$ws.Range("E5").Value = 528611

# --- Pakistan (row 15) ---
$ws.Range("B15").Value = 257914
$ws.Range("C15").Value = 2145
$ws.Range("D15").Value = 178737
$ws.Range("E15").Value = 73751
$ws.Range("G15").Value = 40
$ws.Range("H15").Value = 5426

# --- Kazajistan (row 33) ---
$ws.Range("B33").Value = 65188
$ws.Range("C33").Value = 1674
$ws.Range("D33").Value = 39066
$ws.Range("E33").Value = 25747

# --- Honduras (row 55) ---
$ws.Range("B55").Value = 30036
$ws.Range("C55").Value = 930
$ws.Range("D55").Value = 3379
$ws.Range("E55").Value = 25832
$ws.Range("G55").Value = 18
$ws.Range("H55").Value = 825

# --- Kirguistan (row 72) ---
$ws.Range("B72").Value = 12282
$ws.Range("C72").Value = 305
$ws.Range("D72").Value = 3712
$ws.Range("E72").Value = 8405
$ws.Range("G72").Value = 5
$ws.Range("H72").Value = 165

# --- Australia (row 74) ---
$ws.Range("B74").Value = 10810
$ws.Range("C74").Value = 323
$ws.Range("D74").Value = 8035
$ws.Range("E74").Value = 2662

# --- Reunion (row 156) ---
$ws.Range("B156").Value = 608
$ws.Range("E156").Value = 133

# --- Mauricio (row 163) ---
$ws.Range("B163").Value = 343
$ws.Range("D163").Value = 331

# --- Guyana (row 167) ---
$ws.Range("B167").Value = 313
$ws.Range("E167").Value = 139
$ws.Range("H167").Value = 18

# --- Mongolia (row 169) ---
$ws.Range("D169").Value = 209
$ws.Range("E169").Value = 52

# --- Camboya (row 177) ---
$ws.Range("B177").Value = 166
$ws.Range("C177").Value = 1
$ws.Range("E177").Value = 33
